# "Toasts system & reload module after import"
# Journal de travail: log a new entry for the implementation work done on
# the toast/notification system and the module reload-after-import feature.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Journal de travail")

# Row 52 sits right after the last logged entry (row 51) and before the
# blank spacer rows that precede the totals row. Clone the date-formatted
# style (numFmtId 14, the same style already used by column A's other
# entries) from A51 onto A52 before writing the new values so the new row
# renders as a date like its neighbours.
$ws.Range("A51").Copy()
$ws.Range("A52").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A52").Value = 45085
$ws.Range("B52").Value = "Implémentation"
$ws.Range("C52").Value = 4
$ws.Range("D52").Value = 'Frontend: recherche, ajout via zip, base édition des écrans, gestion des "notifications"'

# C62 = SUM(C2:C61) recalculates automatically to 158.5 once the new 4h
# entry in C52 is included.
